$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.225.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = "'1.673.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.59%  '
$ws.Range("D5").Value = "'211.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.29%  '
$ws.Range("D6").Value = "'0.5282"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.76%  '
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.58%  '
$ws.Range("D8").Value = "'0.2649"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.30%  '
$ws.Range("D9").Value = "'0.06283"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.65%  '
$ws.Range("D10").Value = "'21.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").Value = "'0.07568"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").Value = "'1.676.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").Value = "'4.467"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("D14").Value = "'0.5610"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.29%  '
$ws.Range("D15").Value = "'66.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").Value = "'0.000008018"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.78%  '
$ws.Range("D17").Value = "'26.028.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.75%  '
$ws.Range("D18").Value = "'1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").Value = "'4.820"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.68%  '
$ws.Range("D20").Value = "'187.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.02%  '
$ws.Range("D21").Value = "'10.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.37%  '
$ws.Range("D22").Value = "'6.219"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("D24").Value = "'149.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("D25").Value = "'0.1256"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.86%  '
$ws.Range("D26").Value = "'7.582"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.09%  '
$ws.Range("D27").Value = "'15.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").Value = "'0.06233"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").Value = "'1.361"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("D30").Value = "'1.284"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("D31").Value = "'3.502"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.17%  '
$ws.Range("D32").Value = "'3.432"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.58%  '
$ws.Range("D33").Value = "'1.633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.50%  '
$ws.Range("D34").Value = "'1.003"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.67%  '
$ws.Range("D35").Value = "'0.6055"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("D36").Value = "'2.413"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = "'2.752"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("D38").Value = "'6.114"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").Value = "'0.01621"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("D40").Value = "'1.101.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("D41").Value = "'0.8737"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("D42").Value = "'1.006"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.96%  '
$ws.Range("D43").Value = "'99.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("D44").Value = "'1.825.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.35%  '
$ws.Range("D46").Value = "'56.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.72%  '
$ws.Range("D47").Value = "'1.006"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("D48").Value = "'8.032"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.26%  '
$ws.Range("D49").Value = "'0.05230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("D50").Value = "'0.4256"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("D51").Value = "'5.991"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.31%  '
